# Fixed update to excel issue
#
# 1. Rename the "Requested quantity" headers on the existing sheets.
# 2. Add a new "PO Forecast" worksheet with forecast data (ds, PO_Forecast,
#    yhat_lower, yhat_upper) after "Monthly Trend".

$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# Add the new sheet after the last existing sheet ("Monthly Trend").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Match the look of the other sheets' header row (bold, centered, bordered)
# by copying the existing header format instead of re-building it by hand.
$wsWeekly.Range("A1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122) # xlPasteFormats

# Data rows
$dates = @(44934.99999999999, 44983.99999999999, 44990.99999999999, 45025.99999999999, 45060.99999999999, 45067.99999999999, 45074.99999999999, 45081.99999999999, 45088.99999999999, 45095.99999999999, 45102.99999999999, 45109.99999999999, 45116.99999999999)
$lowers = @(7.999999989681808, 7.999999990229848, 7.999999989468121, 7.999999990098512, 7.999999990148513, 7.999999990199069, 7.999999988403463, 7.999999987262052, 7.999999986434977, 7.999999983361267, 7.999999979641575, 7.999999973230418, 7.99999996748069)
$uppers = @(8.000000009509698, 8.000000009870023, 8.000000009667247, 8.000000009501912, 8.000000010347184, 8.000000010300358, 8.000000010776064, 8.000000012155805, 8.000000013951668, 8.000000017033985, 8.000000021920039, 8.000000026198981, 8.000000033874441)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $wsForecast.Cells.Item($row, 1).Value = $dates[$i]
    $wsForecast.Cells.Item($row, 2).Value = 8
    $wsForecast.Cells.Item($row, 3).Value = $lowers[$i]
    $wsForecast.Cells.Item($row, 4).Value = $uppers[$i]
}

# Match the date formatting used in column A of the other sheets.
$wsWeekly.Range("A2").Copy()
$wsForecast.Range("A2:A14").PasteSpecial(-4122) # xlPasteFormats

# Keep the originally active sheet selected, since the source diff does not
# touch the workbook's active-tab setting.
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
